$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "add more keyword in no.4" -- append two more keyword rows under the
# existing "สี่" (no.4) tag group, right after the last one (row 81),
# reusing the two previously-blank rows 82/83.
$ws.Range("A82").Value = "สี่"
$ws.Range("B82").Value = "แรก"
$ws.Range("A83").Value = "สี่"
$ws.Range("B83").Value = "ขั้นแรก"

# Normalize row heights across the used area (and the stray fully-blank
# rows below it) so every row falls back to the sheet's default height
# instead of carrying an explicit per-row height -- this also drops the
# now-pointless fully-blank row records from the sheet.
$ws.Range("A1:B224").EntireRow.AutoFit()

# Leave the selection where the user would have ended up after typing
# the new keyword pair.
[void]$ws.Range("C84").Select()
